$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 2 body first: "Flags Ops" (D2:E6) ---
$ws.Range("D2").Value = "Positif"
$ws.Range("E2").Value = 2
$ws.Range("D3").Value = "Négatif"
$ws.Range("E3").Value = 0
$ws.Range("D4").Value = "Rapproché"
$ws.Range("E4").Value = 1
$ws.Range("D5").Value = "Rappel"
$ws.Range("E5").Value = 32
$ws.Range("D6").Value = "Ventilé"
$ws.Range("E6").Value = 256

# --- Retitle table 1 header: "Flags opérations" (A1:B1) ---
$ws.Range("A1").Value = "Flags opérations"
$ws.Range("B1").Value = "Valeur"

# --- Table 1 body (unchanged values, B2:B13) ---
$ws.Range("A2").Value = "Somme positive"
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = "Somme négative"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "Somme positive Rapprochée"
$ws.Range("B4").Value = 3
$ws.Range("A5").Value = "Somme négative Rapprochée"
$ws.Range("B5").Value = 1
$ws.Range("A6").Value = "Somme positive Rappel"
$ws.Range("B6").Value = 34
$ws.Range("A7").Value = "Somme négative Rappel"
$ws.Range("B7").Value = 32
$ws.Range("A8").Value = "Somme positive Ventilée"
$ws.Range("B8").Value = 258
$ws.Range("A9").Value = "Somme négative Ventilée"
$ws.Range("B9").Value = 256
$ws.Range("A10").Value = "Somme positive Ventilée Rapprochée"
$ws.Range("B10").Value = 259
$ws.Range("A11").Value = "Somme négative Ventilée Rapprochée"
$ws.Range("B11").Value = 257
$ws.Range("A12").Value = "Somme positive Ventilée Rappel"
$ws.Range("B12").Value = 290
$ws.Range("A13").Value = "Somme négative Ventilée Rappel"
$ws.Range("B13").Value = 288

# --- Table 3 header + first rows: "Flags Catégories" (A15:B18) ---
$ws.Range("A15").Value = "Flags Catégories"
$ws.Range("B15").Value = "Valeur"
$ws.Range("A16").Value = "Dépense"
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = "Revenu"
$ws.Range("B17").Value = 2
$ws.Range("A18").Value = "Possède un parent"
$ws.Range("B18").Value = 1

# --- Back up to add table 2 header: "Flags Ops" (D1:E1) ---
$ws.Range("D1").Value = "Flags Ops"
$ws.Range("E1").Value = "Valeur"

# --- Finish table 3 remaining rows (A19:B20) ---
$ws.Range("A19").Value = "Possède une Op planifiée"
$ws.Range("B19").Value = 4
$ws.Range("A20").Value = "Possède un budget"
$ws.Range("B20").Value = 8

# --- View / selection state ---
$ws.Range("B16").Select()
$excel.ActiveWindow.ScrollRow = 13
